# Applies the "output_offset" update: extends the data table on "Sheet"
# from columns C:L (10 freq points) out to C:T (18 freq points), with
# refreshed frequency / input-offset / output-offset sample values, and
# moves the active selection to A3. Also nudges the workbook window size.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")
$ws.Activate()

# Row 1 = frequency (Hz): 1 GHz .. 18 GHz
$row1Values = @(1000000000, 2000000000, 3000000000, 4000000000, 5000000000, 6000000000, 7000000000, 8000000000, 9000000000, 10000000000, 11000000000, 12000000000, 13000000000, 14000000000, 15000000000, 16000000000, 17000000000, 18000000000)
# Row 2 = input offset (dB)
$row2Values = @(-20.5, -20.7, -20.7, -20.7, -20.7, -20.8, -20.7, -20.7, -20.8, -20.7, -20.6, -20.8, -20.8, -20.7, -20.6, -20.5, -20.4, -20.3)
# Row 3 = output offset (dB)
$row3Values = @(-20.52, -20.71, -20.7, -20.69, -20.68, -20.77, -20.74, -20.7, -20.82, -20.65, -20.64, -20.77, -20.8, -20.75, -20.62, -20.51, -20.4, -20.35)

for ($i = 0; $i -lt $row1Values.Length; $i++) {
    $col = 3 + $i   # column C = 3
    $ws.Cells.Item(1, $col).Value = $row1Values[$i]
    $ws.Cells.Item(2, $col).Value = $row2Values[$i]
    $ws.Cells.Item(3, $col).Value = $row3Values[$i]
}

# Move / record the active selection on A3.
$ws.Range("A3").Select()

# Resize / reposition the workbook window to match the saved state.
$win = $excel.ActiveWindow
$win.Left = -44475
$win.Top = 5310
$win.Width = 18660
$win.Height = 7530
